$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove Param3 header and column E entirely
$ws.Columns.Item(5).Delete() | Out-Null

# Delete the now-extra row 101 (101 rows -> 100 rows)
$ws.Rows.Item(101).Delete() | Out-Null

$ws.Cells.Item(2,1).Value = "Zelda--param1-00.79--1-40.dac"
$ws.Cells.Item(2,2).Value = "Zelda"
$ws.Cells.Item(2,3).Value = 0.79
$ws.Cells.Item(2,4).Value = 40

$ws.Cells.Item(3,1).Value = "Zelda--param1-00.02--5-14.dac"
$ws.Cells.Item(3,2).Value = "Zelda"
$ws.Cells.Item(3,3).Value = 0.02
$ws.Cells.Item(3,4).Value = 14

$ws.Cells.Item(4,1).Value = "Zelda--param1-00.65--3-20.dac"
$ws.Cells.Item(4,2).Value = "Zelda"
$ws.Cells.Item(4,3).Value = 0.65
$ws.Cells.Item(4,4).Value = 20

$ws.Cells.Item(5,1).Value = "Zelda--param1-00.03--3-19.dac"
$ws.Cells.Item(5,2).Value = "Zelda"
$ws.Cells.Item(5,3).Value = 0.03
$ws.Cells.Item(5,4).Value = 19

$ws.Cells.Item(6,1).Value = "Zelda--param1-00.51--6-05.dac"
$ws.Cells.Item(6,2).Value = "Zelda"
$ws.Cells.Item(6,3).Value = 0.51
$ws.Cells.Item(6,4).Value = 5

$ws.Cells.Item(7,1).Value = "Zelda--param1-00.25--6-12.dac"
$ws.Cells.Item(7,2).Value = "Zelda"
$ws.Cells.Item(7,3).Value = 0.25
$ws.Cells.Item(7,4).Value = 12

$ws.Cells.Item(8,1).Value = "Zelda--param1-00.27--3-29.dac"
$ws.Cells.Item(8,2).Value = "Zelda"
$ws.Cells.Item(8,3).Value = 0.27
$ws.Cells.Item(8,4).Value = 29

$ws.Cells.Item(9,1).Value = "Zelda--param1-00.87--1-31.dac"
$ws.Cells.Item(9,2).Value = "Zelda"
$ws.Cells.Item(9,3).Value = 0.87
$ws.Cells.Item(9,4).Value = 31

$ws.Cells.Item(10,1).Value = "Zelda--param1-00.30--4-14.dac"
$ws.Cells.Item(10,2).Value = "Zelda"
$ws.Cells.Item(10,3).Value = 0.3
$ws.Cells.Item(10,4).Value = 14

$ws.Cells.Item(11,1).Value = "Zelda--param1-00.15--2-10.dac"
$ws.Cells.Item(11,2).Value = "Zelda"
$ws.Cells.Item(11,3).Value = 0.15
$ws.Cells.Item(11,4).Value = 10

$ws.Cells.Item(12,1).Value = "Zelda--param1-00.77--4-23.dac"
$ws.Cells.Item(12,2).Value = "Zelda"
$ws.Cells.Item(12,3).Value = 0.77
$ws.Cells.Item(12,4).Value = 23

$ws.Cells.Item(13,1).Value = "Zelda--param1-00.06--1-04.dac"
$ws.Cells.Item(13,2).Value = "Zelda"
$ws.Cells.Item(13,3).Value = 0.06
$ws.Cells.Item(13,4).Value = 4

$ws.Cells.Item(14,1).Value = "Zelda--param1-00.83--6-11.dac"
$ws.Cells.Item(14,2).Value = "Zelda"
$ws.Cells.Item(14,3).Value = 0.83
$ws.Cells.Item(14,4).Value = 11

$ws.Cells.Item(15,1).Value = "Zelda--param1-00.35.dac"
$ws.Cells.Item(15,2).Value = "Zelda"
$ws.Cells.Item(15,3).Value = 0.35
$ws.Cells.Item(15,4).Value = ""

$ws.Cells.Item(16,1).Value = "Zelda--param1-00.39--3-06.dac"
$ws.Cells.Item(16,2).Value = "Zelda"
$ws.Cells.Item(16,3).Value = 0.39
$ws.Cells.Item(16,4).Value = 6

$ws.Cells.Item(17,1).Value = "Zelda--param1-00.22--1-32.dac"
$ws.Cells.Item(17,2).Value = "Zelda"
$ws.Cells.Item(17,3).Value = 0.22
$ws.Cells.Item(17,4).Value = 32

$ws.Cells.Item(18,1).Value = "Zelda--param1-00.83--2-35.dac"
$ws.Cells.Item(18,2).Value = "Zelda"
$ws.Cells.Item(18,3).Value = 0.83
$ws.Cells.Item(18,4).Value = 35

$ws.Cells.Item(19,1).Value = "Zelda--param1-00.95--2-26.dac"
$ws.Cells.Item(19,2).Value = "Zelda"
$ws.Cells.Item(19,3).Value = 0.95
$ws.Cells.Item(19,4).Value = 26

$ws.Cells.Item(20,1).Value = "Zelda--param1-00.52--3-04.dac"
$ws.Cells.Item(20,2).Value = "Zelda"
$ws.Cells.Item(20,3).Value = 0.52
$ws.Cells.Item(20,4).Value = 4

$ws.Cells.Item(21,1).Value = "Zelda--param1-00.72--4-18.dac"
$ws.Cells.Item(21,2).Value = "Zelda"
$ws.Cells.Item(21,3).Value = 0.72
$ws.Cells.Item(21,4).Value = 18

$ws.Cells.Item(22,1).Value = "lofi--param1-00.74.dac"
$ws.Cells.Item(22,2).Value = "lofi"
$ws.Cells.Item(22,3).Value = 0.74
$ws.Cells.Item(22,4).Value = ""

$ws.Cells.Item(23,1).Value = "lofi--param1-00.54.dac"
$ws.Cells.Item(23,2).Value = "lofi"
$ws.Cells.Item(23,3).Value = 0.54
$ws.Cells.Item(23,4).Value = ""

$ws.Cells.Item(24,1).Value = "lofi--param1-00.38.dac"
$ws.Cells.Item(24,2).Value = "lofi"
$ws.Cells.Item(24,3).Value = 0.38
$ws.Cells.Item(24,4).Value = ""

$ws.Cells.Item(25,1).Value = "lofi--param1-00.60--Cymatics.dac"
$ws.Cells.Item(25,2).Value = "lofi"
$ws.Cells.Item(25,3).Value = 0.6
$ws.Cells.Item(25,4).Value = ""

$ws.Cells.Item(26,1).Value = "lofi--param1-00.91.dac"
$ws.Cells.Item(26,2).Value = "lofi"
$ws.Cells.Item(26,3).Value = 0.91
$ws.Cells.Item(26,4).Value = ""

$ws.Cells.Item(27,1).Value = "lofi--param1-00.62--Cymatics.dac"
$ws.Cells.Item(27,2).Value = "lofi"
$ws.Cells.Item(27,3).Value = 0.62
$ws.Cells.Item(27,4).Value = ""

$ws.Cells.Item(28,1).Value = "lofi--param1-00.93.dac"
$ws.Cells.Item(28,2).Value = "lofi"
$ws.Cells.Item(28,3).Value = 0.93
$ws.Cells.Item(28,4).Value = ""

$ws.Cells.Item(29,1).Value = "lofi--param1-00.25--Cymatics.dac"
$ws.Cells.Item(29,2).Value = "lofi"
$ws.Cells.Item(29,3).Value = 0.25
$ws.Cells.Item(29,4).Value = ""

$ws.Cells.Item(30,1).Value = "lofi--param1-00.22.dac"
$ws.Cells.Item(30,2).Value = "lofi"
$ws.Cells.Item(30,3).Value = 0.22
$ws.Cells.Item(30,4).Value = ""

$ws.Cells.Item(31,1).Value = "lofi--param1-00.14.dac"
$ws.Cells.Item(31,2).Value = "lofi"
$ws.Cells.Item(31,3).Value = 0.14
$ws.Cells.Item(31,4).Value = ""

$ws.Cells.Item(32,1).Value = "lofi--param1-00.06.dac"
$ws.Cells.Item(32,2).Value = "lofi"
$ws.Cells.Item(32,3).Value = 0.06
$ws.Cells.Item(32,4).Value = ""

$ws.Cells.Item(33,1).Value = "lofi--param1-00.92.dac"
$ws.Cells.Item(33,2).Value = "lofi"
$ws.Cells.Item(33,3).Value = 0.92
$ws.Cells.Item(33,4).Value = ""

$ws.Cells.Item(34,1).Value = "lofi--param1-00.50.dac"
$ws.Cells.Item(34,2).Value = "lofi"
$ws.Cells.Item(34,3).Value = 0.5
$ws.Cells.Item(34,4).Value = ""

$ws.Cells.Item(35,1).Value = "lofi--param1-00.02.dac"
$ws.Cells.Item(35,2).Value = "lofi"
$ws.Cells.Item(35,3).Value = 0.02
$ws.Cells.Item(35,4).Value = ""

$ws.Cells.Item(36,1).Value = "lofi--param1-00.01.dac"
$ws.Cells.Item(36,2).Value = "lofi"
$ws.Cells.Item(36,3).Value = 0.01
$ws.Cells.Item(36,4).Value = ""

$ws.Cells.Item(37,1).Value = "lofi--param1-00.20.dac"
$ws.Cells.Item(37,2).Value = "lofi"
$ws.Cells.Item(37,3).Value = 0.2
$ws.Cells.Item(37,4).Value = ""

$ws.Cells.Item(38,1).Value = "lofi--param1-01.00.dac"
$ws.Cells.Item(38,2).Value = "lofi"
$ws.Cells.Item(38,3).Value = 1.0
$ws.Cells.Item(38,4).Value = ""

$ws.Cells.Item(39,1).Value = "lofi--param1-00.18.dac"
$ws.Cells.Item(39,2).Value = "lofi"
$ws.Cells.Item(39,3).Value = 0.18
$ws.Cells.Item(39,4).Value = ""

$ws.Cells.Item(40,1).Value = "lofi--param1-00.35.dac"
$ws.Cells.Item(40,2).Value = "lofi"
$ws.Cells.Item(40,3).Value = 0.35
$ws.Cells.Item(40,4).Value = ""

$ws.Cells.Item(41,1).Value = "Fusion--param1-00.25.dac"
$ws.Cells.Item(41,2).Value = "Fusion"
$ws.Cells.Item(41,3).Value = 0.25
$ws.Cells.Item(41,4).Value = ""

$ws.Cells.Item(42,1).Value = "Fusion--param1-00.65.dac"
$ws.Cells.Item(42,2).Value = "Fusion"
$ws.Cells.Item(42,3).Value = 0.65
$ws.Cells.Item(42,4).Value = ""

$ws.Cells.Item(43,1).Value = "Fusion--param1-00.40.dac"
$ws.Cells.Item(43,2).Value = "Fusion"
$ws.Cells.Item(43,3).Value = 0.4
$ws.Cells.Item(43,4).Value = ""

$ws.Cells.Item(44,1).Value = "Fusion--param1-00.22.dac"
$ws.Cells.Item(44,2).Value = "Fusion"
$ws.Cells.Item(44,3).Value = 0.22
$ws.Cells.Item(44,4).Value = ""

$ws.Cells.Item(45,1).Value = "Fusion--param1-00.09.dac"
$ws.Cells.Item(45,2).Value = "Fusion"
$ws.Cells.Item(45,3).Value = 0.09
$ws.Cells.Item(45,4).Value = ""

$ws.Cells.Item(46,1).Value = "Fusion--param1-00.96--segment_03..dac"
$ws.Cells.Item(46,2).Value = "Fusion"
$ws.Cells.Item(46,3).Value = 0.96
$ws.Cells.Item(46,4).Value = ""

$ws.Cells.Item(47,1).Value = "Fusion--param1-00.24.dac"
$ws.Cells.Item(47,2).Value = "Fusion"
$ws.Cells.Item(47,3).Value = 0.24
$ws.Cells.Item(47,4).Value = ""

$ws.Cells.Item(48,1).Value = "Fusion--param1-00.20.dac"
$ws.Cells.Item(48,2).Value = "Fusion"
$ws.Cells.Item(48,3).Value = 0.2
$ws.Cells.Item(48,4).Value = ""

$ws.Cells.Item(49,1).Value = "Fusion--param1-00.51.dac"
$ws.Cells.Item(49,2).Value = "Fusion"
$ws.Cells.Item(49,3).Value = 0.51
$ws.Cells.Item(49,4).Value = ""

$ws.Cells.Item(50,1).Value = "Fusion--param1-00.27.dac"
$ws.Cells.Item(50,2).Value = "Fusion"
$ws.Cells.Item(50,3).Value = 0.27
$ws.Cells.Item(50,4).Value = ""

$ws.Cells.Item(51,1).Value = "Fusion--param1-00.78.dac"
$ws.Cells.Item(51,2).Value = "Fusion"
$ws.Cells.Item(51,3).Value = 0.78
$ws.Cells.Item(51,4).Value = ""

$ws.Cells.Item(52,1).Value = "Fusion--param1-00.11.dac"
$ws.Cells.Item(52,2).Value = "Fusion"
$ws.Cells.Item(52,3).Value = 0.11
$ws.Cells.Item(52,4).Value = ""

$ws.Cells.Item(53,1).Value = "Fusion--param1-00.45.dac"
$ws.Cells.Item(53,2).Value = "Fusion"
$ws.Cells.Item(53,3).Value = 0.45
$ws.Cells.Item(53,4).Value = ""

$ws.Cells.Item(54,1).Value = "Fusion--param1-00.37.dac"
$ws.Cells.Item(54,2).Value = "Fusion"
$ws.Cells.Item(54,3).Value = 0.37
$ws.Cells.Item(54,4).Value = ""

$ws.Cells.Item(55,1).Value = "Fusion--param1-00.92--segment_09..dac"
$ws.Cells.Item(55,2).Value = "Fusion"
$ws.Cells.Item(55,3).Value = 0.92
$ws.Cells.Item(55,4).Value = ""

$ws.Cells.Item(56,1).Value = "Fusion--param1-00.79.dac"
$ws.Cells.Item(56,2).Value = "Fusion"
$ws.Cells.Item(56,3).Value = 0.79
$ws.Cells.Item(56,4).Value = ""

$ws.Cells.Item(57,1).Value = "Fusion--param1-00.34--segment_03..dac"
$ws.Cells.Item(57,2).Value = "Fusion"
$ws.Cells.Item(57,3).Value = 0.34
$ws.Cells.Item(57,4).Value = ""

$ws.Cells.Item(58,1).Value = "Fusion--param1-00.42.dac"
$ws.Cells.Item(58,2).Value = "Fusion"
$ws.Cells.Item(58,3).Value = 0.42
$ws.Cells.Item(58,4).Value = ""

$ws.Cells.Item(59,1).Value = "Fusion--param1-00.53.dac"
$ws.Cells.Item(59,2).Value = "Fusion"
$ws.Cells.Item(59,3).Value = 0.53
$ws.Cells.Item(59,4).Value = ""

$ws.Cells.Item(60,1).Value = "Fusion--param1-00.82.dac"
$ws.Cells.Item(60,2).Value = "Fusion"
$ws.Cells.Item(60,3).Value = 0.82
$ws.Cells.Item(60,4).Value = ""

$ws.Cells.Item(61,1).Value = "8bit--param1-00.19--40.dac"
$ws.Cells.Item(61,2).Value = "8bit"
$ws.Cells.Item(61,3).Value = 0.19
$ws.Cells.Item(61,4).Value = ""

$ws.Cells.Item(62,1).Value = "8bit--param1-00.62--23.dac"
$ws.Cells.Item(62,2).Value = "8bit"
$ws.Cells.Item(62,3).Value = 0.62
$ws.Cells.Item(62,4).Value = ""

$ws.Cells.Item(63,1).Value = "8bit--param1-00.81.dac"
$ws.Cells.Item(63,2).Value = "8bit"
$ws.Cells.Item(63,3).Value = 0.81
$ws.Cells.Item(63,4).Value = ""

$ws.Cells.Item(64,1).Value = "8bit--param1-00.01.dac"
$ws.Cells.Item(64,2).Value = "8bit"
$ws.Cells.Item(64,3).Value = 0.01
$ws.Cells.Item(64,4).Value = ""

$ws.Cells.Item(65,1).Value = "8bit--param1-00.50.dac"
$ws.Cells.Item(65,2).Value = "8bit"
$ws.Cells.Item(65,3).Value = 0.5
$ws.Cells.Item(65,4).Value = ""

$ws.Cells.Item(66,1).Value = "8bit--param1-00.55.dac"
$ws.Cells.Item(66,2).Value = "8bit"
$ws.Cells.Item(66,3).Value = 0.55
$ws.Cells.Item(66,4).Value = ""

$ws.Cells.Item(67,1).Value = "8bit--param1-00.85.dac"
$ws.Cells.Item(67,2).Value = "8bit"
$ws.Cells.Item(67,3).Value = 0.85
$ws.Cells.Item(67,4).Value = ""

$ws.Cells.Item(68,1).Value = "8bit--param1-00.71--31.dac"
$ws.Cells.Item(68,2).Value = "8bit"
$ws.Cells.Item(68,3).Value = 0.71
$ws.Cells.Item(68,4).Value = ""

$ws.Cells.Item(69,1).Value = "8bit--param1-00.15--43.dac"
$ws.Cells.Item(69,2).Value = "8bit"
$ws.Cells.Item(69,3).Value = 0.15
$ws.Cells.Item(69,4).Value = ""

$ws.Cells.Item(70,1).Value = "8bit--param1-00.60.dac"
$ws.Cells.Item(70,2).Value = "8bit"
$ws.Cells.Item(70,3).Value = 0.6
$ws.Cells.Item(70,4).Value = ""

$ws.Cells.Item(71,1).Value = "8bit--param1-00.06.dac"
$ws.Cells.Item(71,2).Value = "8bit"
$ws.Cells.Item(71,3).Value = 0.06
$ws.Cells.Item(71,4).Value = ""

$ws.Cells.Item(72,1).Value = "8bit--param1-00.14.dac"
$ws.Cells.Item(72,2).Value = "8bit"
$ws.Cells.Item(72,3).Value = 0.14
$ws.Cells.Item(72,4).Value = ""

$ws.Cells.Item(73,1).Value = "8bit--param1-00.62--01.dac"
$ws.Cells.Item(73,2).Value = "8bit"
$ws.Cells.Item(73,3).Value = 0.62
$ws.Cells.Item(73,4).Value = ""

$ws.Cells.Item(74,1).Value = "8bit--param1-00.88.dac"
$ws.Cells.Item(74,2).Value = "8bit"
$ws.Cells.Item(74,3).Value = 0.88
$ws.Cells.Item(74,4).Value = ""

$ws.Cells.Item(75,1).Value = "8bit--param1-00.82.dac"
$ws.Cells.Item(75,2).Value = "8bit"
$ws.Cells.Item(75,3).Value = 0.82
$ws.Cells.Item(75,4).Value = ""

$ws.Cells.Item(76,1).Value = "8bit--param1-00.67--19.dac"
$ws.Cells.Item(76,2).Value = "8bit"
$ws.Cells.Item(76,3).Value = 0.67
$ws.Cells.Item(76,4).Value = ""

$ws.Cells.Item(77,1).Value = "8bit--param1-00.49.dac"
$ws.Cells.Item(77,2).Value = "8bit"
$ws.Cells.Item(77,3).Value = 0.49
$ws.Cells.Item(77,4).Value = ""

$ws.Cells.Item(78,1).Value = "8bit--param1-00.40.dac"
$ws.Cells.Item(78,2).Value = "8bit"
$ws.Cells.Item(78,3).Value = 0.4
$ws.Cells.Item(78,4).Value = ""

$ws.Cells.Item(79,1).Value = "8bit--param1-00.71--46.dac"
$ws.Cells.Item(79,2).Value = "8bit"
$ws.Cells.Item(79,3).Value = 0.71
$ws.Cells.Item(79,4).Value = ""

$ws.Cells.Item(80,1).Value = "8bit--param1-00.12.dac"
$ws.Cells.Item(80,2).Value = "8bit"
$ws.Cells.Item(80,3).Value = 0.12
$ws.Cells.Item(80,4).Value = ""

$ws.Cells.Item(81,1).Value = "duduk--param1-00.22.dac"
$ws.Cells.Item(81,2).Value = "duduk"
$ws.Cells.Item(81,3).Value = 0.22
$ws.Cells.Item(81,4).Value = ""

$ws.Cells.Item(82,1).Value = "duduk--param1-00.51.dac"
$ws.Cells.Item(82,2).Value = "duduk"
$ws.Cells.Item(82,3).Value = 0.51
$ws.Cells.Item(82,4).Value = ""

$ws.Cells.Item(83,1).Value = "duduk--param1-00.77.dac"
$ws.Cells.Item(83,2).Value = "duduk"
$ws.Cells.Item(83,3).Value = 0.77
$ws.Cells.Item(83,4).Value = ""

$ws.Cells.Item(84,1).Value = "duduk--param1-00.85.dac"
$ws.Cells.Item(84,2).Value = "duduk"
$ws.Cells.Item(84,3).Value = 0.85
$ws.Cells.Item(84,4).Value = ""

$ws.Cells.Item(85,1).Value = "duduk--param1-00.70.dac"
$ws.Cells.Item(85,2).Value = "duduk"
$ws.Cells.Item(85,3).Value = 0.7
$ws.Cells.Item(85,4).Value = ""

$ws.Cells.Item(86,1).Value = "duduk--param1-00.11.dac"
$ws.Cells.Item(86,2).Value = "duduk"
$ws.Cells.Item(86,3).Value = 0.11
$ws.Cells.Item(86,4).Value = ""

$ws.Cells.Item(87,1).Value = "duduk--param1-00.38.dac"
$ws.Cells.Item(87,2).Value = "duduk"
$ws.Cells.Item(87,3).Value = 0.38
$ws.Cells.Item(87,4).Value = ""

$ws.Cells.Item(88,1).Value = "duduk--param1-00.55.dac"
$ws.Cells.Item(88,2).Value = "duduk"
$ws.Cells.Item(88,3).Value = 0.55
$ws.Cells.Item(88,4).Value = ""

$ws.Cells.Item(89,1).Value = "duduk--param1-00.66.dac"
$ws.Cells.Item(89,2).Value = "duduk"
$ws.Cells.Item(89,3).Value = 0.66
$ws.Cells.Item(89,4).Value = ""

$ws.Cells.Item(90,1).Value = "duduk--param1-00.88.dac"
$ws.Cells.Item(90,2).Value = "duduk"
$ws.Cells.Item(90,3).Value = 0.88
$ws.Cells.Item(90,4).Value = ""

$ws.Cells.Item(91,1).Value = "duduk--param1-00.02.dac"
$ws.Cells.Item(91,2).Value = "duduk"
$ws.Cells.Item(91,3).Value = 0.02
$ws.Cells.Item(91,4).Value = ""

$ws.Cells.Item(92,1).Value = "duduk--param1-00.80.dac"
$ws.Cells.Item(92,2).Value = "duduk"
$ws.Cells.Item(92,3).Value = 0.8
$ws.Cells.Item(92,4).Value = ""

$ws.Cells.Item(93,1).Value = "duduk--param1-00.54--KSHMR_Duduk_14_One_Shot_Fm.dac"
$ws.Cells.Item(93,2).Value = "duduk"
$ws.Cells.Item(93,3).Value = 0.54
$ws.Cells.Item(93,4).Value = ""

$ws.Cells.Item(94,1).Value = "duduk--param1-00.26.dac"
$ws.Cells.Item(94,2).Value = "duduk"
$ws.Cells.Item(94,3).Value = 0.26
$ws.Cells.Item(94,4).Value = ""

$ws.Cells.Item(95,1).Value = "duduk--param1-00.58.dac"
$ws.Cells.Item(95,2).Value = "duduk"
$ws.Cells.Item(95,3).Value = 0.58
$ws.Cells.Item(95,4).Value = ""

$ws.Cells.Item(96,1).Value = "duduk--param1-00.63.dac"
$ws.Cells.Item(96,2).Value = "duduk"
$ws.Cells.Item(96,3).Value = 0.63
$ws.Cells.Item(96,4).Value = ""

$ws.Cells.Item(97,1).Value = "duduk--param1-00.23.dac"
$ws.Cells.Item(97,2).Value = "duduk"
$ws.Cells.Item(97,3).Value = 0.23
$ws.Cells.Item(97,4).Value = ""

$ws.Cells.Item(98,1).Value = "duduk--param1-00.74.dac"
$ws.Cells.Item(98,2).Value = "duduk"
$ws.Cells.Item(98,3).Value = 0.74
$ws.Cells.Item(98,4).Value = ""

$ws.Cells.Item(99,1).Value = "duduk--param1-00.93.dac"
$ws.Cells.Item(99,2).Value = "duduk"
$ws.Cells.Item(99,3).Value = 0.93
$ws.Cells.Item(99,4).Value = ""

$ws.Cells.Item(100,1).Value = "duduk--param1-00.76.dac"
$ws.Cells.Item(100,2).Value = "duduk"
$ws.Cells.Item(100,3).Value = 0.76
$ws.Cells.Item(100,4).Value = ""

Write-Host "done"